{"js": "// Update the worksheet date and all two-digit x two-digit multiplication\n// problems/answers to the new set of values.\nconst replacements = [\n  [\"2025-08-08 Friday\", \"2025-08-09 Saturday\"],\n  [\"83\u00d748=3984\", \"64\u00d711=704\"],\n  [\"97\u00d753=5141\", \"63\u00d723=1449\"],\n  [\"73\u00d751=3723\", \"38\u00d794=3572\"],\n  [\"78\u00d789=6942\", \"55\u00d771=3905\"],\n  [\"91\u00d735=3185\", \"91\u00d789=8099\"],\n  [\"18\u00d711=198\", \"15\u00d726=390\"],\n  [\"96\u00d788=8448\", \"70\u00d768=4760\"],\n  [\"61\u00d771=4331\", \"65\u00d765=4225\"],\n  [\"57\u00d784=4788\", \"63\u00d746=2898\"],\n  [\"47\u00d734=1598\", \"83\u00d774=6142\"],\n  [\"39\u00d758=2262\", \"80\u00d774=5920\"],\n  [\"81\u00d740=3240\", \"84\u00d742=3528\"],\n  [\"59\u00d725=1475\", \"31\u00d765=2015\"],\n  [\"28\u00d779=2212\", \"30\u00d718=540\"],\n  [\"32\u00d785=2720\", \"67\u00d736=2412\"],\n  [\"80\u00d796=7680\", \"40\u00d753=2120\"],\n  [\"29\u00d732=928\", \"86\u00d726=2236\"],\n  [\"88\u00d738=3344\", \"89\u00d782=7298\"],\n  [\"89\u00d793=8277\", \"73\u00d724=1752\"],\n  [\"57\u00d730=1710\", \"41\u00d760=2460\"],\n  [\"65\u00d799=6435\", \"54\u00d792=4968\"],\n  [\"91\u00d795=8645\", \"60\u00d782=4920\"],\n  [\"40\u00d715=600\", \"90\u00d754=4860\"],\n  [\"97\u00d794=9118\", \"32\u00d786=2752\"],\n  [\"75\u00d767=5025\", \"86\u00d716=1376\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all two-digit x two-digit multiplication\n# problems/answers to the new set of values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-08-08 Friday\"; New = \"2025-08-09 Saturday\" },\n    @{ Old = \"83\u00d748=3984\"; New = \"64\u00d711=704\" },\n    @{ Old = \"97\u00d753=5141\"; New = \"63\u00d723=1449\" },\n    @{ Old = \"73\u00d751=3723\"; New = \"38\u00d794=3572\" },\n    @{ Old = \"78\u00d789=6942\"; New = \"55\u00d771=3905\" },\n    @{ Old = \"91\u00d735=3185\"; New = \"91\u00d789=8099\" },\n    @{ Old = \"18\u00d711=198\"; New = \"15\u00d726=390\" },\n    @{ Old = \"96\u00d788=8448\"; New = \"70\u00d768=4760\" },\n    @{ Old = \"61\u00d771=4331\"; New = \"65\u00d765=4225\" },\n    @{ Old = \"57\u00d784=4788\"; New = \"63\u00d746=2898\" },\n    @{ Old = \"47\u00d734=1598\"; New = \"83\u00d774=6142\" },\n    @{ Old = \"39\u00d758=2262\"; New = \"80\u00d774=5920\" },\n    @{ Old = \"81\u00d740=3240\"; New = \"84\u00d742=3528\" },\n    @{ Old = \"59\u00d725=1475\"; New = \"31\u00d765=2015\" },\n    @{ Old = \"28\u00d779=2212\"; New = \"30\u00d718=540\" },\n    @{ Old = \"32\u00d785=2720\"; New = \"67\u00d736=2412\" },\n    @{ Old = \"80\u00d796=7680\"; New = \"40\u00d753=2120\" },\n    @{ Old = \"29\u00d732=928\"; New = \"86\u00d726=2236\" },\n    @{ Old = \"88\u00d738=3344\"; New = \"89\u00d782=7298\" },\n    @{ Old = \"89\u00d793=8277\"; New = \"73\u00d724=1752\" },\n    @{ Old = \"57\u00d730=1710\"; New = \"41\u00d760=2460\" },\n    @{ Old = \"65\u00d799=6435\"; New = \"54\u00d792=4968\" },\n    @{ Old = \"91\u00d795=8645\"; New = \"60\u00d782=4920\" },\n    @{ Old = \"40\u00d715=600\"; New = \"90\u00d754=4860\" },\n    @{ Old = \"97\u00d794=9118\"; New = \"32\u00d786=2752\" },\n    @{ Old = \"75\u00d767=5025\"; New = \"86\u00d716=1376\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n\n$d.Save()\n"}
